$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header: extend the sequence with P1=14, Q1=15, reusing O1's style (bold/border)
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$excel.CutCopyMode = $false

# Swap the I and K column values for rows 2-25
$iVals = $ws.Range("I2:I25").Value2
$kVals = $ws.Range("K2:K25").Value2
$ws.Range("I2:I25").Value2 = $kVals
$ws.Range("K2:K25").Value2 = $iVals

# Swap the M and O column values for rows 2-25
$mVals = $ws.Range("M2:M25").Value2
$oVals = $ws.Range("O2:O25").Value2
$ws.Range("M2:M25").Value2 = $oVals
$ws.Range("O2:O25").Value2 = $mVals

# Add new data columns P and Q for rows 2-25, filled with 2
$ws.Range("P2:P25").Value2 = 2
$ws.Range("Q2:Q25").Value2 = 2
